$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instructions")

# Part B: load/store, jump, upper
# lb/lh/lw (I*-type) opcode bit change: row 8, columns J:L (1 -> 0)
$ws.Range("J8:L8").Value = 0

# jalr (J-type block) ImmSel_2 bit change: N24 (1 -> 0)
$ws.Range("N24").Value = 0

# Add missing "I-type" header label above the jalr column (N19)
$ws.Range("N19").Value = "I-type"
$ws.Range("N19").Font.Bold = $true

# Update the active selection to reflect where the author ended up
$ws.Range("O20").Select()
